$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Append the new daily-push row (row 60) with the same layout as the
# existing data rows: date text, weekday text, hour number, rank number.
$row = 60
$ws.Cells.Item($row, 1).NumberFormat = "@"
$ws.Cells.Item($row, 1).Value = "2025/10/04"
$ws.Cells.Item($row, 1).Style = "Normal"
$ws.Cells.Item($row, 2).Value = "土"
$ws.Cells.Item($row, 3).Value = 20
$ws.Cells.Item($row, 4).Value = 42
